$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.856.94'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '3.419.18'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '570.47'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '157.96'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.422.55'
$ws.Range("E8").Value = '  +0.09%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.569'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -8.97%  '
$ws.Range("E10").Value = '  +1.62%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.98%  '
$ws.Range("E12").Value = '  -4.05%  '
$ws.Range("D13").Value = '4.007.39'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("E16").Value = '  -7.83%  '
$ws.Range("D17").Value = '63.942.37'
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("D18").Value = '3.389.72'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("E19").Value = '  -3.54%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.62'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.86%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '382.10'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.39%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.81'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  -5.60%  '
$ws.Range("E26").Value = '  -3.18%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.68'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.11%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  -1.85%  '
$ws.Range("E31").Value = '  -5.67%  '
$ws.Range("E32").Value = '  -0.69%  '
$ws.Range("E33").Value = '  +0.05%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '22.90'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("E36").Value = '  -5.76%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '160.80'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("E38").Value = '  +10.04%  '
$ws.Range("E39").Value = '  -3.89%  '
$ws.Range("D40").Value = '2.810.39'
$ws.Range("E40").Value = '  -2.21%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '25.96'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.88%  '
$ws.Range("E42").Value = '  -5.04%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '43.06'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.49%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '26.34'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  -7.86%  '
$ws.Range("E46").Value = '  -5.55%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0304'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.68%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +9.17%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '333.68'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.46%  '
$ws.Range("E50").Value = '  -3.45%  '
$ws.Range("E51").Value = '  -5.26%  '
